$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 8000363
$ws.Range("I33").Value = 4762301
$ws.Range("J33").Value = 25000188
$ws.Range("K33").Value = 4762301
$ws.Range("L33").Value = 25000188
$ws.Range("M33").Value = -4762072
$ws.Range("N33").Value = -25000646

$ws.Range("H96").Value = 6671885.5
$ws.Range("I96").Value = 7586
$ws.Range("K96").Value = 22758
$ws.Range("M96").Value = -21385

$ws.Range("H106").Value = 11500
$ws.Range("I106").Value = 3000
$ws.Range("K106").Value = 3000
$ws.Range("M106").Value = -2369

$ws.Range("H113").Value = 168494.5
$ws.Range("I113").Value = 2158.3333
$ws.Range("J113").Value = 334830.66
$ws.Range("K113").Value = 2158.3333
$ws.Range("L113").Value = 334830.66
$ws.Range("M113").Value = 1095.6667
$ws.Range("N113").Value = -341338.66

$ws.Range("H116").Value = 6917.6665
$ws.Range("I116").Value = 6638.5
$ws.Range("J116").Value = 7141
$ws.Range("K116").Value = 6638.5
$ws.Range("L116").Value = 7141
$ws.Range("M116").Value = -3196.5
$ws.Range("N116").Value = -14025

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 13148000
$ws.Range("I2").Value = 13148000
$ws.Range("K2").Value = 13148000
$ws.Range("M2").Value = -13147887

$ws.Range("H30").Value = 1201.25
$ws.Range("I30").Value = 1567.6666
$ws.Range("J30").Value = 102
$ws.Range("K30").Value = 1567.6666
$ws.Range("L30").Value = 102
$ws.Range("M30").Value = -1417.6666
$ws.Range("N30").Value = -402

$ws.Range("H32").Value = 27271.512
$ws.Range("I32").Value = 27902.262
$ws.Range("K32").Value = 27902.262
$ws.Range("M32").Value = -27615.262

$ws.Range("H61").Value = 3458.7144
$ws.Range("I61").Value = 2938.6667
$ws.Range("J61").Value = 17500
$ws.Range("K61").Value = 2938.6667
$ws.Range("L61").Value = 17500
$ws.Range("M61").Value = -2726.6667
$ws.Range("N61").Value = -17924

$ws.Range("H63").Value = 9895.450000000001
$ws.Range("I63").Value = 9000
$ws.Range("J63").Value = 9942.579
$ws.Range("K63").Value = 9000
$ws.Range("L63").Value = 9942.579
$ws.Range("M63").Value = -8314
$ws.Range("N63").Value = -11314.579

$ws.Range("H66").Value = 9895.450000000001
$ws.Range("I66").Value = 9000
$ws.Range("J66").Value = 9942.579
$ws.Range("K66").Value = 45000
$ws.Range("L66").Value = 49712.895
$ws.Range("M66").Value = -41568
$ws.Range("N66").Value = -56576.895

$ws.Range("H74").Value = 591991.3
$ws.Range("I74").Value = 1429301.6
$ws.Range("K74").Value = 1429301.6
$ws.Range("M74").Value = -1428427.6

$ws.Range("H77").Value = 591991.3
$ws.Range("I77").Value = 1429301.6
$ws.Range("K77").Value = 7146508
$ws.Range("M77").Value = -7142140

$ws.Range("H110").Value = 15001017
$ws.Range("I110").Value = 17308754
$ws.Range("K110").Value = 17308754
$ws.Range("M110").Value = -17306709

$ws.Range("H116").Value = 13148000
$ws.Range("I116").Value = 13148000
$ws.Range("K116").Value = 13148000
$ws.Range("M116").Value = -13145706

$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").ClearContents()

$ws.Range("H132").Value = 4954.9316
$ws.Range("I132").Value = 2806.0833
$ws.Range("J132").Value = 14624.75
$ws.Range("K132").Value = 8418.249899999999
$ws.Range("L132").Value = 43874.25
$ws.Range("M132").Value = -5888.249899999999
$ws.Range("N132").Value = -48934.25

$ws.Range("H136").Value = 3458.7144
$ws.Range("I136").Value = 2938.6667
$ws.Range("J136").Value = 17500
$ws.Range("K136").Value = 8816.000100000001
$ws.Range("L136").Value = 52500
$ws.Range("M136").Value = -6266.000100000001
$ws.Range("N136").Value = -57600

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 13148000
$ws.Range("I3").Value = 13148000
$ws.Range("K3").Value = 13148000
$ws.Range("M3").Value = -13147886

$ws.Range("H11").Value = 219.125
$ws.Range("I11").Value = 236.14285
$ws.Range("J11").Value = 100
$ws.Range("K11").Value = 236.14285
$ws.Range("L11").Value = 100
$ws.Range("M11").Value = -96.14285000000001
$ws.Range("N11").Value = -380

$ws.Range("H37").Value = 817
$ws.Range("I37").Value = 817
$ws.Range("K37").Value = 817
$ws.Range("M37").Value = -680

$ws.Range("H94").Value = 1408.65
$ws.Range("I94").Value = 898.5
$ws.Range("K94").Value = 898.5
$ws.Range("M94").Value = -447.5

$ws.Range("H105").Value = 200051470
$ws.Range("I105").Value = 200051470
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 200051470
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = -200049723
$ws.Range("N105").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H10").Value = 1990.8334
$ws.Range("I10").Value = 2907.5
$ws.Range("J10").Value = 157.5
$ws.Range("K10").Value = 2907.5
$ws.Range("L10").Value = 157.5
$ws.Range("M10").Value = -2768.5
$ws.Range("N10").Value = -435.5

$ws.Range("H122").Value = 50927.95
$ws.Range("I122").Value = 59762.94
$ws.Range("J122").Value = 863
$ws.Range("K122").Value = 179288.82
$ws.Range("L122").Value = 2589
$ws.Range("M122").Value = -176838.82
$ws.Range("N122").Value = -7489

$ws.Range("H132").Value = 205179.17
$ws.Range("I132").Value = 3643.75
$ws.Range("K132").Value = 10931.25
$ws.Range("M132").Value = -8401.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H47").Value = 143461.28
$ws.Range("I47").Value = 200143.8
$ws.Range("K47").Value = 600431.3999999999
$ws.Range("M47").Value = -600000.3999999999

$ws.Range("H74").Value = 10335.5
$ws.Range("I74").Value = 9006.5
$ws.Range("J74").Value = 11000
$ws.Range("K74").Value = 27019.5
$ws.Range("L74").Value = 33000
$ws.Range("M74").Value = -25958.5
$ws.Range("N74").Value = -35122

$ws.Range("H77").Value = 10335.5
$ws.Range("I77").Value = 9006.5
$ws.Range("J77").Value = 11000
$ws.Range("K77").Value = 81058.5
$ws.Range("L77").Value = 99000
$ws.Range("M77").Value = -75754.5
$ws.Range("N77").Value = -109608

$ws.Range("H80").Value = 4066
$ws.Range("I80").Value = 3998
$ws.Range("J80").Value = 4100
$ws.Range("K80").Value = 11994
$ws.Range("L80").Value = 12300
$ws.Range("M80").Value = -11058
$ws.Range("N80").Value = -14172

$ws.Range("H83").Value = 4066
$ws.Range("I83").Value = 3998
$ws.Range("J83").Value = 4100
$ws.Range("K83").Value = 35982
$ws.Range("L83").Value = 36900
$ws.Range("M83").Value = -31302
$ws.Range("N83").Value = -46260

$ws.Range("H97").Value = 506.41666
$ws.Range("I97").Value = 516.5
$ws.Range("J97").Value = 496.33334
$ws.Range("K97").Value = 1549.5
$ws.Range("L97").Value = 1489.00002
$ws.Range("M97").Value = -1053.5
$ws.Range("N97").Value = -2481.00002

$ws.Range("H98").Value = 295
$ws.Range("I98").Value = 295
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 885
$ws.Range("L98").Value = 0
$ws.Range("M98").Value = 613
$ws.Range("N98").ClearContents()

$ws.Range("H132").Value = 1891.6
$ws.Range("J132").Value = 2665.1667
$ws.Range("L132").Value = 23986.5003
$ws.Range("N132").Value = -29046.5003

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1294.5
$ws.Range("I97").Value = 1294.5
$ws.Range("K97").Value = 1294.5
$ws.Range("M97").Value = -798.5

$ws.Range("H107").Value = 1217
$ws.Range("J107").Value = 1367.6666
$ws.Range("L107").Value = 1367.6666
$ws.Range("N107").Value = -5207.6666

$ws.Range("H122").Value = 7758.4443
$ws.Range("J122").Value = 4949
$ws.Range("L122").Value = 14847
$ws.Range("N122").Value = -19747

$ws.Range("H132").Value = 4537.1304
$ws.Range("J132").Value = 10332.833
$ws.Range("L132").Value = 30998.499
$ws.Range("N132").Value = -36058.499

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3198.9092
$ws.Range("J22").Value = 4796.643
$ws.Range("L22").Value = 4796.643
$ws.Range("N22").Value = -5386.643

$ws.Range("H27").Value = 3198.9092
$ws.Range("J27").Value = 4796.643
$ws.Range("L27").Value = 4796.643
$ws.Range("N27").Value = -5010.643

$ws.Range("H46").Value = 2945805.8
$ws.Range("J46").Value = 6227
$ws.Range("L46").Value = 6227
$ws.Range("N46").Value = -6603

$ws.Range("H93").Value = 2296.6875
$ws.Range("I93").Value = 1448.2778
$ws.Range("K93").Value = 1448.2778
$ws.Range("M93").Value = -200.2778000000001

$ws.Range("H100").Value = 13894915
$ws.Range("I100").Value = 22731224
$ws.Range("K100").Value = 22731224
$ws.Range("M100").Value = -22730683

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2687.2856
$ws.Range("I122").Value = 2559.3684
$ws.Range("J122").Value = 3902.5
$ws.Range("K122").Value = 7678.1052
$ws.Range("L122").Value = 11707.5
$ws.Range("M122").Value = -5228.1052
$ws.Range("N122").Value = -16607.5

$ws.Range("H132").Value = 6320.25
$ws.Range("I132").Value = 2889.7778
$ws.Range("J132").Value = 8669.3334
$ws.Range("M132").Value = -6139.3334
